$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Font.Name = "Noto Sans"
Write-Output "done"
